$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-11 (Adam2-Itga9 LR pairs); only the
# columns affected by the TPM recompute are touched, everything else is
# left exactly as it was.
$updates = @{
    2 = @{ "I"=0.6113651253405055; "J"=0.6113651253405055; "M"=5.854382333333334; "N"=17.563147; "O"=0.1730451459016118; "P"=0.1730451459016118; "Q"=0.1729989494107778; "R"=1.556990544697; "S"=0.105793767313705; "T"=0.105793767313705 }
    3 = @{ "I"=0.6113651253405055; "J"=0.6113651253405055; "O"=0.1208497063316524; "P"=0.1208497063316525; "S"=0.07388329585881397; "T"=0.07388329585881398 }
    4 = @{ "I"=0.6113651253405055; "J"=0.6113651253405055; "M"=11.64342866666667; "N"=34.930286; "O"=0.3441590756630932; "P"=0.3441590756630932; "Q"=0.3440671982428888; "R"=3.096604784185999; "S"=0.2104068564298395; "T"=0.2104068564298395 }
    5 = @{ "I"=0.6113651253405055; "J"=0.6113651253405055; "M"=0.1645376666666667; "N"=0.493613; "O"=0.004863441250245888; "P"=0.004863441250245888; "Q"=0.004862142895888888; "R"=0.043759286063; "S"=0.002973338369542762; "T"=0.002973338369542762 }
    6 = @{ "I"=0.6113651253405055; "J"=0.6113651253405055; "M"=12.08065233333333; "N"=36.241957; "O"=0.3570826308533967; "P"=0.3570826308533967; "Q"=0.3569873033341111; "R"=3.212885730007; "S"=0.2183078673686043; "T"=0.2183078673686043 }
    7 = @{ "E"=1; "F"=0.3333333333333333; "G"=0.01878466666666667; "H"=0.056354; "I"=0.3886348746594945; "J"=0.3886348746594945; "M"=5.854382333333334; "N"=17.563147; "O"=0.1730451459016118; "P"=0.1730451459016118; "Q"=0.1099726206708889; "R"=0.989753586038; "S"=0.06725137858790685; "T"=0.06725137858790685 }
    8 = @{ "E"=1; "F"=0.3333333333333333; "G"=0.01878466666666667; "H"=0.056354; "I"=0.3886348746594945; "J"=0.3886348746594945; "O"=0.1208497063316524; "P"=0.1208497063316525; "Q"=0.07680168572977777; "R"=0.691215171568; "S"=0.04696641047283847; "T"=0.04696641047283848 }
    9 = @{ "E"=1; "F"=0.3333333333333333; "G"=0.01878466666666667; "H"=0.056354; "I"=0.3886348746594945; "J"=0.3886348746594945; "M"=11.64342866666667; "N"=34.930286; "O"=0.3441590756630932; "P"=0.3441590756630932; "Q"=0.2187179263604444; "R"=1.968461337244; "S"=0.1337522192332537; "T"=0.1337522192332537 }
    10 = @{ "E"=1; "F"=0.3333333333333333; "G"=0.01878466666666667; "H"=0.056354; "I"=0.3886348746594945; "J"=0.3886348746594945; "M"=0.1645376666666667; "N"=0.493613; "O"=0.004863441250245888; "P"=0.004863441250245888; "Q"=0.003090785222444445; "R"=0.027817067002; "S"=0.001890102880703126; "T"=0.001890102880703126 }
    11 = @{ "E"=1; "F"=0.3333333333333333; "G"=0.01878466666666667; "H"=0.056354; "I"=0.3886348746594945; "J"=0.3886348746594945; "M"=12.08065233333333; "N"=36.241957; "O"=0.3570826308533967; "P"=0.3570826308533967; "Q"=0.2269310271975556; "R"=2.042379244778; "S"=0.1387747634847924; "T"=0.1387747634847924 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
